$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1323.9584
$ws.Range("J32").Value = 1393.7
$ws.Range("L32").Value = 1393.7
$ws.Range("N32").Value = -2045.7
$ws.Range("H55").Value = 820.25
$ws.Range("I55").Value = 993.6667
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 993.6667
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = -779.6667
$ws.Range("N55").Value = -728
$ws.Range("H113").Value = 2915.5
$ws.Range("I113").Value = 2801
$ws.Range("J113").Value = 3030
$ws.Range("K113").Value = 2801
$ws.Range("L113").Value = 3030
$ws.Range("M113").Value = 453
$ws.Range("N113").Value = -9538
$ws.Range("H116").Value = 5159.6
$ws.Range("I116").Value = 5331.52
$ws.Range("J116").Value = 4300
$ws.Range("K116").Value = 5331.52
$ws.Range("L116").Value = 4300
$ws.Range("M116").Value = -1889.52
$ws.Range("N116").Value = -11184

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3183.382
$ws.Range("I32").Value = 1931.6385
$ws.Range("J32").Value = 20499.166
$ws.Range("K32").Value = 1931.6385
$ws.Range("L32").Value = 20499.166
$ws.Range("M32").Value = -1644.6385
$ws.Range("N32").Value = -21073.166
$ws.Range("H76").Value = 31500
$ws.Range("J76").Value = 31500
$ws.Range("L76").Value = 31500
$ws.Range("N76").Value = -32176
$ws.Range("H79").Value = 31500
$ws.Range("J79").Value = 31500
$ws.Range("L79").Value = 31500
$ws.Range("N79").Value = -33840
$ws.Range("H82").Value = 29980
$ws.Range("J82").Value = 29980
$ws.Range("L82").Value = 29980
$ws.Range("N82").Value = -30702
$ws.Range("H85").Value = 29980
$ws.Range("J85").Value = 29980
$ws.Range("L85").Value = 29980
$ws.Range("N85").Value = -32476
$ws.Range("H92").Value = 28909.6
$ws.Range("J92").Value = 28909.6
$ws.Range("L92").Value = 28909.6
$ws.Range("N92").Value = -33901.6
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H101").Value = 59999.5
$ws.Range("J101").Value = 59999.5
$ws.Range("L101").Value = 59999.5
$ws.Range("N101").Value = -66489.5
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H119").Value = 65000
$ws.Range("J119").Value = 65000
$ws.Range("L119").Value = 65000
$ws.Range("N119").Value = -74676
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 2800.805
$ws.Range("I132").Value = 1117.5172
$ws.Range("J132").Value = 6868.75
$ws.Range("K132").Value = 3352.5516
$ws.Range("L132").Value = 20606.25
$ws.Range("M132").Value = -822.5515999999998
$ws.Range("N132").Value = -25666.25
$ws.Range("H135").Value = 67999.75
$ws.Range("J135").Value = 67999.75
$ws.Range("L135").Value = 67999.75
$ws.Range("N135").Value = -78139.75
$ws.Range("H139").Value = 43369
$ws.Range("J139").Value = 43369
$ws.Range("L139").Value = 43369
$ws.Range("N139").Value = -53649

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13137.429
$ws.Range("I20").Value = 1153.5385
$ws.Range("J20").Value = 32611.25
$ws.Range("K20").Value = 1153.5385
$ws.Range("L20").Value = 32611.25
$ws.Range("M20").Value = -906.5385000000001
$ws.Range("N20").Value = -33105.25
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 1260.7693
$ws.Range("I86").Value = 957.7143
$ws.Range("J86").Value = 1614.3334
$ws.Range("K86").Value = 957.7143
$ws.Range("L86").Value = 1614.3334
$ws.Range("M86").Value = 165.2857
$ws.Range("N86").Value = -3860.3334
$ws.Range("H89").Value = 1260.7693
$ws.Range("I89").Value = 957.7143
$ws.Range("J89").Value = 1614.3334
$ws.Range("K89").Value = 4788.5715
$ws.Range("L89").Value = 8071.666999999999
$ws.Range("M89").Value = 827.4285
$ws.Range("N89").Value = -19303.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 42307784
$ws.Range("I12").Value = 100000100
$ws.Range("J12").Value = 89.73333
$ws.Range("K12").Value = 300000300
$ws.Range("L12").Value = 269.19999
$ws.Range("M12").Value = -300000127
$ws.Range("N12").Value = -615.19999
$ws.Range("H131").Value = 1887720.1
$ws.Range("I131").Value = 9091343
$ws.Range("J131").Value = 1057.1428
$ws.Range("K131").Value = 27274029
$ws.Range("L131").Value = 3171.4284
$ws.Range("M131").Value = -27268989
$ws.Range("N131").Value = -13251.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 40001264
$ws.Range("I113").Value = 76924110
$ws.Range("J113").Value = 1511.25
$ws.Range("K113").Value = 76924110
$ws.Range("L113").Value = 1511.25
$ws.Range("M113").Value = -76921940
$ws.Range("N113").Value = -5851.25
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H123").Value = 20871.834
$ws.Range("J123").Value = 20871.834
$ws.Range("L123").Value = 20871.834
$ws.Range("N123").Value = -25771.834
$ws.Range("H126").Value = 5226.1562
$ws.Range("I126").Value = 7808.5
$ws.Range("K126").Value = 23425.5
$ws.Range("M126").Value = -20955.5
$ws.Range("H132").Value = 4182.5386
$ws.Range("I132").Value = 5404.6
$ws.Range("J132").Value = 3418.75
$ws.Range("K132").Value = 16213.8
$ws.Range("L132").Value = 10256.25
$ws.Range("M132").Value = -13683.8
$ws.Range("N132").Value = -15316.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 85000
$ws.Range("J69").Value = 80000
$ws.Range("L69").Value = 80000
$ws.Range("N69").Value = -81622
$ws.Range("H72").Value = 85000
$ws.Range("J72").Value = 80000
$ws.Range("L72").Value = 240000
$ws.Range("N72").Value = -248112
$ws.Range("H104").Value = 14375
$ws.Range("J104").Value = 14375
$ws.Range("L104").Value = 14375
$ws.Range("N104").Value = -21363
$ws.Range("H110").Value = 39900
$ws.Range("J110").Value = 39900
$ws.Range("L110").Value = 39900
$ws.Range("N110").Value = -48080
$ws.Range("H119").Value = 28500
$ws.Range("J119").Value = 28500
$ws.Range("L119").Value = 28500
$ws.Range("N119").Value = -38176
$ws.Range("H132").Value = 12066544
$ws.Range("I132").Value = 18339634
$ws.Range("J132").Value = 2909.7693
$ws.Range("K132").Value = 55018902
$ws.Range("L132").Value = 8729.3079
$ws.Range("M132").Value = -55016372
$ws.Range("N132").Value = -13789.3079
$ws.Range("H133").Value = 40312
$ws.Range("J133").Value = 40312
$ws.Range("L133").Value = 40312
$ws.Range("N133").Value = -45372

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 8853
$ws.Range("J41").Value = 8853
$ws.Range("L41").Value = 8853
$ws.Range("N41").Value = -9633
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
